$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D ("Tipo"), shifting it (and its data) to column E
$ws.Columns.Item(4).Insert()

# Header + value for the new "MAE" column, matching the header style used
# by the other column headers (bold, bordered, centered)
$ws.Range("D1").Value = "MAE"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D2").Value = 0.1911874935925034

# Recomputed MSE value for the existing row
$ws.Range("B2").Value = 0.07796894984218575
